$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Cell D48 (FraxShare) needs to stay text ("8.60") so the trailing zero survives;
# force Text format just for that one cell before writing the value.
$ws.Range("D48").NumberFormat = "@"

$ws.Range("D2").Value = '43.729.53'
$ws.Range("E2").Value = '  -0.20%  '

$ws.Range("D3").Value = '2.290.12'
$ws.Range("E3").Value = '  -0.15%  '

$ws.Range("E4").Value = '  -0.09%  '

$ws.Range("D5").Value = '115.22'
$ws.Range("E5").Value = '  +1.84%  '

$ws.Range("D6").Value = '266.13'

$ws.Range("D7").Value = '0.645'
$ws.Range("E7").Value = '  +2.95%  '

$ws.Range("E8").Value = '  +0.11%  '

$ws.Range("E9").Value = '  -1.15%  '

$ws.Range("D10").Value = '47.43'
$ws.Range("E10").Value = '  -1.35%  '

$ws.Range("E11").Value = '  -1.05%  '

$ws.Range("D12").Value = '9.16'
$ws.Range("E12").Value = '  -0.03%  '

$ws.Range("E13").Value = '  +1.73%  '

$ws.Range("D14").Value = '15.41'
$ws.Range("E14").Value = '  -2.38%  '

$ws.Range("D15").Value = '2.635.02'
$ws.Range("E15").Value = '  -0.12%  '

$ws.Range("D16").Value = '0.874'
$ws.Range("E16").Value = '  +2.58%  '

$ws.Range("D17").Value = '2.291.61'
$ws.Range("E17").Value = '  +0.05%  '

$ws.Range("D18").Value = '43.644.07'
$ws.Range("E18").Value = '  -0.11%  '

$ws.Range("E19").Value = '  +0.08%  '

$ws.Range("E20").Value = '  +0.23%  '

$ws.Range("D21").Value = '72.51'
$ws.Range("E21").Value = '  +0.30%  '

$ws.Range("D22").Value = '2.45'
$ws.Range("E22").Value = '  +0.05%  '

$ws.Range("D23").Value = '236.69'
$ws.Range("E23").Value = '  +1.78%  '

$ws.Range("B24").Value = 'InternetComputer(DFINITY)'
$ws.Range("C24").Value = 'https://coinranking.com/coin/aMNLwaUbY+internetcomputerdfinity-icp'
$ws.Range("D24").Value = '9.48'
$ws.Range("E24").Value = '  -3.33%  '

$ws.Range("B25").Value = 'PancakeSwap'
$ws.Range("C25").Value = 'https://coinranking.com/coin/ncYFcP709+pancakeswap-cake'
$ws.Range("D25").Value = '2.88'
$ws.Range("E25").Value = '  +1.57%  '

$ws.Range("E26").Value = '  +1.89%  '

$ws.Range("D27").Value = '11.58'
$ws.Range("E27").Value = '  -0.76%  '

$ws.Range("D28").Value = '41.61'
$ws.Range("E28").Value = '  -1.16%  '

$ws.Range("E29").Value = '  -0.31%  '

$ws.Range("E30").Value = '  -0.73%  '

$ws.Range("D31").Value = '173.94'
$ws.Range("E31").Value = '  -0.83%  '

$ws.Range("D32").Value = '21.76'
$ws.Range("E32").Value = '  +1.06%  '

$ws.Range("E33").Value = '  -1.68%  '

$ws.Range("D34").Value = '5.68'
$ws.Range("E34").Value = '  -0.24%  '

$ws.Range("E35").Value = '  +2.29%  '

$ws.Range("E36").Value = '  +5.38%  '

$ws.Range("D37").Value = '4.69'
$ws.Range("E37").Value = '  +0.53%  '

$ws.Range("E38").Value = '  +3.89%  '

$ws.Range("E39").Value = '  -1.72%  '

$ws.Range("E40").Value = '  +7.50%  '

$ws.Range("D41").Value = '14.37'
$ws.Range("E41").Value = '  +3.92%  '

$ws.Range("D42").Value = '74.33'
$ws.Range("E42").Value = '  +0.20%  '

$ws.Range("D43").Value = '0.236'
$ws.Range("E43").Value = '  -2.72%  '

$ws.Range("D44").Value = '6.01'
$ws.Range("E44").Value = '  -4.82%  '

$ws.Range("E45").Value = '  -0.05%  '

$ws.Range("D46").Value = '1.37'
$ws.Range("E46").Value = '  -0.99%  '

$ws.Range("D47").Value = '1.28'
$ws.Range("E47").Value = '  +3.60%  '

$ws.Range("D48").Value = '8.60'
$ws.Range("E48").Value = '  -2.15%  '

$ws.Range("D49").Value = '73.32'
$ws.Range("E49").Value = '  +34.37%  '

$ws.Range("E50").Value = '  +0.41%  '

$ws.Range("D51").Value = '100.69'
$ws.Range("E51").Value = '  -2.19%  '
